# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (values that Excel will not auto-convert to a number,
# e.g. multi-dot "thousand-separated" prices and padded percentage strings).
$ws.Range("D2").Value = '58.405.65'
$ws.Range("E2").Value = '  -3.73%  '
$ws.Range("D3").Value = '2.697.51'
$ws.Range("E3").Value = '  -7.03%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").Value = '  -5.66%  '
$ws.Range("E6").Value = '  -3.06%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  -5.58%  '
$ws.Range("D9").Value = '2.708.01'
$ws.Range("E9").Value = '  -6.89%  '
$ws.Range("E10").Value = '  +1.14%  '
$ws.Range("E11").Value = '  -3.61%  '
$ws.Range("E12").Value = '  -4.81%  '
$ws.Range("E13").Value = '  +0.83%  '
$ws.Range("D14").Value = '3.171.92'
$ws.Range("E14").Value = '  -6.97%  '
$ws.Range("D15").Value = '58.498.53'
$ws.Range("E15").Value = '  -3.51%  '
$ws.Range("E16").Value = '  -6.11%  '
$ws.Range("D17").Value = '2.705.56'
$ws.Range("E17").Value = '  -6.89%  '
$ws.Range("E18").Value = '  -6.11%  '
$ws.Range("E19").Value = '  -6.68%  '
$ws.Range("E20").Value = '  -6.99%  '
$ws.Range("E21").Value = '  -8.73%  '
$ws.Range("E22").Value = '  -7.06%  '
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("E25").Value = '  -2.54%  '
$ws.Range("E26").Value = '  -6.88%  '
$ws.Range("E27").Value = '  -4.73%  '
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("E29").Value = '  -5.38%  '
$ws.Range("D30").Value = '0.0₃0817'
$ws.Range("E30").Value = '  -5.67%  '
$ws.Range("E32").Value = '  -3.52%  '
$ws.Range("E33").Value = '  -5.40%  '
$ws.Range("E34").Value = '  -0.34%  '
$ws.Range("E35").Value = '  -4.31%  '
$ws.Range("E36").Value = '  -5.40%  '
$ws.Range("E37").Value = '  -7.01%  '
$ws.Range("E38").Value = '  -8.32%  '
$ws.Range("E39").Value = '  -6.51%  '
$ws.Range("E40").Value = '  -4.35%  '
$ws.Range("D41").Value = '2.176.03'
$ws.Range("E41").Value = '  -6.32%  '
$ws.Range("E42").Value = '  -8.92%  '
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("E44").Value = '  -4.75%  '
$ws.Range("E45").Value = '  -8.06%  '
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("E47").Value = '  -9.90%  '
$ws.Range("E48").Value = '  -4.85%  '
$ws.Range("E49").Value = '  -7.02%  '
$ws.Range("E50").Value = '  -5.64%  '
$ws.Range("E51").Value = '  -4.13%  '

# Numeric-looking text values (e.g. "62.80") must be forced to stay text,
# otherwise Excel auto-converts them to a number and drops the trailing zero.
# A leading apostrophe forces text entry; ClearFormats() then strips the
# resulting quote-prefix cell styling so the cell keeps its original (default)
# appearance while remaining a text value.
$ws.Range("D5").Value = "'500.03"
$ws.Range("D6").Value = "'139.14"
$ws.Range("D8").Value = "'0.526"
$ws.Range("D21").Value = "'333.18"
$ws.Range("D23").Value = "'0.997"
$ws.Range("D25").Value = "'62.80"
$ws.Range("D28").Value = "'0.996"
$ws.Range("D29").Value = "'7.40"
$ws.Range("D32").Value = "'19.04"
$ws.Range("D34").Value = "'150.70"
$ws.Range("D35").Value = "'5.36"
$ws.Range("D37").Value = "'0.928"
$ws.Range("D39").Value = "'35.24"
$ws.Range("D42").Value = "'1.37"
$ws.Range("D44").Value = "'0.0554"
$ws.Range("D47").Value = "'18.66"
$ws.Range("D51").Value = "'17.76"

$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D51").ClearFormats()
